$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I ("Year") to make room for "Percentage"
$ws.Columns("I").Insert()

# New header + values for the inserted "Percentage" column
$ws.Range("I1").Value = "Percentage"
$ws.Range("I2").Value = 8
$ws.Range("I3").Value = 3.6

# Re-fit the surrounding column widths (mirrors Excel's own best-fit pass
# after inserting a column / typing new header text)
$ws.Columns("B").ColumnWidth = 18.666666666666668
$ws.Columns("C").ColumnWidth = 18.330729166666668
$ws.Columns("D").ColumnWidth = 17.498697916666668
$ws.Columns("E").ColumnWidth = 14.830729166666666
$ws.Columns("F").ColumnWidth = 9.666666666666666
$ws.Columns("G").ColumnWidth = 9.330729166666666
$ws.Columns("I").ColumnWidth = 17.498697916666668

# Keep selection on I2 as in the saved file
$ws.Range("I2").Select()
